# Applies a weekly re-shuffle of the price-record rows (2-21) in the
# "Hortaliza, Agrícola del Norte S.A. de Arica - Repollo" sheet.
#
# Only the per-record columns Fecha (D), Calidad (I), Volumen (J),
# Precio mínimo (K), Precio máximo (L), Precio promedio ponderado (M)
# and Precio $/Kg (P) are shuffled between rows; all the other columns
# (A, B, C, E, F, G, H, N, O, Q, R) stay identical since they carry the
# same constant values on every row. Capture the original values first
# (a simple row permutation) and then write them back in their new rows
# so that source and destination overlaps do not clobber each other.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target row -> Source row (i.e. row $target gets the data that used to
# live in row $source before the edit).
$rowMap = @{
    2  = 12
    3  = 5
    4  = 15
    5  = 16
    6  = 17
    7  = 3
    8  = 4
    9  = 9
    10 = 10
    11 = 13
    12 = 19
    13 = 18
    14 = 11
    15 = 20
    16 = 6
    17 = 7
    18 = 2
    19 = 8
    20 = 21
    21 = 14
}

$cols = @("D", "I", "J", "K", "L", "M", "P")

# Snapshot the original values for the columns we are about to shuffle.
# NOTE: use Value2 (not Value) to read back raw/unformatted values -
# Value2 returns numbers as numbers (e.g. the date serial) and strings
# as strings, which is what we need to write back unchanged.
$snapshot = @{}
foreach ($row in 2..21) {
    $rowValues = @{}
    foreach ($col in $cols) {
        $rowValues[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowValues
}

# Write the shuffled values back into place.
foreach ($targetRow in 2..21) {
    $sourceRow = $rowMap[$targetRow]
    $sourceValues = $snapshot[$sourceRow]
    foreach ($col in $cols) {
        $ws.Range("$col$targetRow").Value = $sourceValues[$col]
    }
}
